# Applies the commit's changes to the "empleados" sheet:
#  - Adds column Q formulas: INT(YEARFRAC(C,TODAY()))  (rows 4-32)
#  - Adds column R formulas: IF(O<5000,300,IF(AND(O>=5000,O<=10000),IF(ISNUMBER(F),200,100),0))  (rows 4-32)
#  - Moves the active selection to R8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 gets its own (non-shared) formulas, matching the original file's
# pattern where row 4 holds the "template" formula and rows 5:32 share it.
$ws.Range("Q4").Formula = "=INT(YEARFRAC(C4,TODAY()))"
$ws.Range("R4").Formula = "=IF(O4<5000,300,IF(AND(O4>=5000,O4<=10000),IF(ISNUMBER(F4),200,100),0))"

# Rows 5:32 share one formula each (Excel auto-detects & stores as shared
# formulas when the same relative formula is applied across a range).
$ws.Range("Q5:Q32").Formula = "=INT(YEARFRAC(C5,TODAY()))"
$ws.Range("R5:R32").Formula = "=IF(O5<5000,300,IF(AND(O5>=5000,O5<=10000),IF(ISNUMBER(F5),200,100),0))"

# Restore focus/selection to R8, as in the edited workbook.
$null = $ws.Range("R8").Select()
